$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row -------------------------------------------------------
# B1 header renamed from "Date" to "Clock.Today" (keeps text format style)
$ws.Range("B1").Value = "Clock.Today"

# --- Column A: was text "Simulation" repeated, now numeric simulation id
$ws.Range("A2").Value = 1
$ws.Range("A3").Value = 2
$ws.Range("A4").Value = 3

# --- New columns E (Wheat.Grain.Wt) and F (ValueWithSpace) headers ----
$ws.Range("E1").NumberFormat = "@"
$ws.Range("E1").Value = "Wheat.Grain.Wt"
$ws.Range("F1").NumberFormat = "@"
$ws.Range("F1").Value = "ValueWithSpace"

# F4 holds a single space value (text) - set before G1 so shared-string
# ordering matches (it must be inserted ahead of "Wheat.Grain.N")
$ws.Range("F4").NumberFormat = "@"
$ws.Range("F4").Value = " "

# --- New column G (Wheat.Grain.N) header -------------------------------
$ws.Range("G1").NumberFormat = "@"
$ws.Range("G1").Value = "Wheat.Grain.N"

# --- Remaining numeric data for E, F, G --------------------------------
$ws.Range("E2").Value = 0

$ws.Range("E3").Value = 10
$ws.Range("F3").Value = 10
$ws.Range("G3").Value = 20

$ws.Range("E4").Value = 1000
$ws.Range("G4").Value = 400

# --- New rows 5 and 6 ----------------------------------------------------
$ws.Range("A5").Value = 1
$ws.Range("B5").NumberFormat = "@"
$ws.Range("B5").Value = "1/01/2000"
$ws.Range("C5").NumberFormat = "@"
$ws.Range("C5").Value = "1 Jan"
$ws.Range("F5").Value = 1000

$ws.Range("A6").Value = 1
$ws.Range("B6").NumberFormat = "@"
$ws.Range("B6").Value = "1/01/2000"
$ws.Range("C6").NumberFormat = "@"
$ws.Range("C6").Value = "1 Jan"

# --- Column widths (closest achievable values given COM width quantization) ---
$ws.Columns.Item(2).ColumnWidth = 16.666666666666668
$ws.Columns.Item(5).ColumnWidth = 17.333333333333332
$ws.Columns.Item(6).ColumnWidth = 15.0

# --- Selection moves to H10 -------------------------------------------
$ws.Range("H10").Select() | Out-Null
